# Coastal Surface Piercing Profilers - Omaha Cal Info update
# Re-point GP05MOAS-GL001 reference designators to GP05MOAS-GL363 and
# switch the active tab from "Moorings" to "Asset_Cal_Info".

$wb = $excel.ActiveWorkbook

# --- Moorings sheet -------------------------------------------------------
$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsMoorings.Range("A2").Value = "GP05MOAS-GL363"

# --- Asset_Cal_Info sheet --------------------------------------------------
$wsAssetCal = $wb.Worksheets.Item("Asset_Cal_Info")
$wsAssetCal.Range("A3").Value = "GP05MOAS-GL363-00-ENG000000"
$wsAssetCal.Range("A4").Value = "GP05MOAS-GL363-01-FLORDM000"
$wsAssetCal.Range("A5").Value = "GP05MOAS-GL363-01-FLORDM000"
$wsAssetCal.Range("A6").Value = "GP05MOAS-GL363-01-FLORDM000"
$wsAssetCal.Range("A7").Value = "GP05MOAS-GL363-01-FLORDM000"
$wsAssetCal.Range("A8").Value = "GP05MOAS-GL363-02-DOSTAM000"
$wsAssetCal.Range("A9").Value = "GP05MOAS-GL363-04-CTDGVM000"

# --- Active tab moves from Moorings to Asset_Cal_Info ----------------------
$wsAssetCal.Activate()
